# Extend the "LevelsSimple" sheet's generated table from W11 (row 80) through
# a new W12 block (rows 81-92), then pad with 8 blank-but-styled rows
# (93-100) so the sheet's used range grows from A1:M80 to A1:M100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone row 80's formatting down through row 100 -------------------
# (gives every new cell the same thin-bottom-border style, s="1", that the
# rest of the generated table uses)
$ws.Range("A80:M80").Copy()
$ws.Range("A81:M100").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Fill in the W12 block (rows 81-92) with the same formula pattern
#        used by every other row of the table ------------------------------
for ($r = 81; $r -le 92; $r++) {
    $p = $r - 1

    $ws.Range("A$r").Formula = "=""W""&B$r&""L""&C$r"
    $ws.Range("B$r").Formula = "=IF(OR(C$p<`$B`$1+1,C$p<B$p),B$p,B$p+1)"
    $ws.Range("C$r").Formula = "=IF(B$p=B$r,C$p+1,1)"
    $ws.Range("E$r").Formula = "=B$r"
    $ws.Range("F$r").Formula = "=E$r+`$B`$1"
    $ws.Range("G$r").Formula = "=C$r"
    $ws.Range("H$r").Formula = "=G$r+`$B`$1"
    $ws.Range("I$r").Formula = "=""+'""&A$r&"" ""&B$r&"" ""&C$r&"" ""&D$r&"" ""&E$r&"" ""&F$r&"" ""&G$r&"" ""&H$r&"" \r\n'"""
    $ws.Range("J$r").Formula = "=E$r*G$r"
    $ws.Range("K$r").Formula = "=E$r*H$r"
    $ws.Range("L$r").Formula = "=F$r*G$r"
    $ws.Range("M$r").Formula = "=F$r*H$r"
}

# Rows 93-100 stay empty (no values/formulas) -- they only picked up the
# copied border styling above, matching the template's trailing blank rows.

# --- 3. Update the view: select I5:I92 (mirrors the sheet's new selection)
#        and scroll so row 4 is at the top of the viewport ------------------
$ws.Range("I5:I92").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
